$wb = $excel.ActiveWorkbook

# Sheet ALC Row 100: Asking for a Friend | Beetle Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 4034.3333
$ws.Cells.Item(100, 10).Value = 4551.5
$ws.Cells.Item(100, 12).Value = 4551.5
$ws.Cells.Item(100, 14).Value = -5633.5

# Sheet ALC Row 107: Another Man's Ink | Enchanted Truegold Ink
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 433.0625
$ws.Cells.Item(107, 10).Value = 448.8
$ws.Cells.Item(107, 12).Value = 448.8
$ws.Cells.Item(107, 14).Value = -4288.8

# Sheet ALC Row 113: Amaro Kart | Starch Glue
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3115.8572
$ws.Cells.Item(113, 9).Value = 2187.7273
$ws.Cells.Item(113, 10).Value = 4136.8
$ws.Cells.Item(113, 11).Value = 2187.7273
$ws.Cells.Item(113, 12).Value = 4136.8
$ws.Cells.Item(113, 13).Value = 1066.2727
$ws.Cells.Item(113, 14).Value = -10644.8

# Sheet ALC Row 116: Growing Up | Growth Formula Kappa
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 14003
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 14003
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).ClearContents()  # L116 removed (was 3730.1538)
$ws.Cells.Item(116, 13).Value = 14003
$ws.Cells.Item(116, 14).Value = -20887

# Sheet ALC Row 129: Practical Command | Commanding Craftsman's Draught
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 278619
$ws.Cells.Item(129, 10).Value = 278619
$ws.Cells.Item(129, 12).Value = 835857
$ws.Cells.Item(129, 14).Value = -845857

# Sheet ALC Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 4349.4
$ws.Cells.Item(132, 9).Value = 4422.8237
$ws.Cells.Item(132, 10).Value = 3933.3333
$ws.Cells.Item(132, 11).Value = 13268.4711
$ws.Cells.Item(132, 12).Value = 11799.9999
$ws.Cells.Item(132, 13).Value = -10738.4711
$ws.Cells.Item(132, 14).Value = -16859.9999

# Sheet ALC Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 84891.664
$ws.Cells.Item(137, 9).Value = 1228.4286
$ws.Cells.Item(137, 10).Value = 202020.2
$ws.Cells.Item(137, 11).Value = 3685.2858
$ws.Cells.Item(137, 12).Value = 606060.6000000001
$ws.Cells.Item(137, 13).Value = -1135.2858
$ws.Cells.Item(137, 14).Value = -611160.6000000001

# Sheet ALC Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1641.9508
$ws.Cells.Item(138, 9).Value = 1112.1945
$ws.Cells.Item(138, 10).Value = 2404.8
$ws.Cells.Item(138, 11).Value = 3336.5835
$ws.Cells.Item(138, 12).Value = 7214.400000000001
$ws.Cells.Item(138, 13).Value = 1803.4165
$ws.Cells.Item(138, 14).Value = -17494.4

# Sheet ARM Row 32: Ingot We Trust | Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 18569.982
$ws.Cells.Item(32, 9).Value = 18872.865
$ws.Cells.Item(32, 11).Value = 18872.865
$ws.Cells.Item(32, 13).Value = -18585.865

# Sheet ARM Row 33: A Leg to Stand On | Heavy Iron Flanchard
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(33, 8).Value = 11675.333
$ws.Cells.Item(33, 9).Value = 6350.6665
$ws.Cells.Item(33, 10).Value = 17000
$ws.Cells.Item(33, 11).Value = 6350.6665
$ws.Cells.Item(33, 12).Value = 17000
$ws.Cells.Item(33, 13).Value = -6021.6665
$ws.Cells.Item(33, 14).Value = -17658

# Sheet ARM Row 45: Hollow Hallmarks | Mythril Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4288.615
$ws.Cells.Item(45, 9).Value = 4318.625
$ws.Cells.Item(45, 10).Value = 4240.6
$ws.Cells.Item(45, 11).Value = 4318.625
$ws.Cells.Item(45, 12).Value = 4240.6
$ws.Cells.Item(45, 13).Value = -3941.625
$ws.Cells.Item(45, 14).Value = -4994.6

# Sheet ARM Row 88: The Mast Chance | Adamantite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 113030.555
$ws.Cells.Item(88, 9).Value = 1799.75
$ws.Cells.Item(88, 10).Value = 202015.2
$ws.Cells.Item(88, 11).Value = 1799.75
$ws.Cells.Item(88, 12).Value = 202015.2
$ws.Cells.Item(88, 13).Value = -1393.75
$ws.Cells.Item(88, 14).Value = -202827.2

# Sheet ARM Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 113030.555
$ws.Cells.Item(91, 9).Value = 1799.75
$ws.Cells.Item(91, 10).Value = 202015.2
$ws.Cells.Item(91, 11).Value = 1799.75
$ws.Cells.Item(91, 12).Value = 202015.2
$ws.Cells.Item(91, 13).Value = -395.75
$ws.Cells.Item(91, 14).Value = -204823.2

# Sheet ARM Row 97: Ore for Me | High Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 1231.6666
$ws.Cells.Item(97, 9).Value = 1542.2222
$ws.Cells.Item(97, 11).Value = 1542.2222
$ws.Cells.Item(97, 13).Value = -1046.2222

# Sheet ARM Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 15769.639
$ws.Cells.Item(132, 9).Value = 1776.409
$ws.Cells.Item(132, 11).Value = 5329.227000000001
$ws.Cells.Item(132, 13).Value = -2799.227000000001

# Sheet ARM Row 135: Forgiveness for My Shins | Ruthenium Sabatons of Fending
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(135, 8).Value = 24740
$ws.Cells.Item(135, 10).Value = 24740
$ws.Cells.Item(135, 12).Value = 24740
$ws.Cells.Item(135, 14).Value = -34880

# Sheet BSM Row 20: Smelt and Dealt | Iron Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3059.9333
$ws.Cells.Item(20, 9).Value = 3223.077
$ws.Cells.Item(20, 10).Value = 1999.5
$ws.Cells.Item(20, 11).Value = 3223.077
$ws.Cells.Item(20, 12).Value = 1999.5
$ws.Cells.Item(20, 13).Value = -2976.077
$ws.Cells.Item(20, 14).Value = -2493.5

# Sheet BSM Row 86: Through Thick and Thin | Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1495.9822
$ws.Cells.Item(86, 9).Value = 1329.8718
$ws.Cells.Item(86, 10).Value = 1877.0588
$ws.Cells.Item(86, 11).Value = 1329.8718
$ws.Cells.Item(86, 12).Value = 1877.0588
$ws.Cells.Item(86, 13).Value = -206.8717999999999
$ws.Cells.Item(86, 14).Value = -4123.0588

# Sheet BSM Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 1495.9822
$ws.Cells.Item(89, 9).Value = 1329.8718
$ws.Cells.Item(89, 10).Value = 1877.0588
$ws.Cells.Item(89, 11).Value = 6649.358999999999
$ws.Cells.Item(89, 12).Value = 9385.294
$ws.Cells.Item(89, 13).Value = -1033.358999999999
$ws.Cells.Item(89, 14).Value = -20617.294

# Sheet BSM Row 94: High Steal | High Steel Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1248.4186
$ws.Cells.Item(94, 9).Value = 1116.4242
$ws.Cells.Item(94, 10).Value = 1684
$ws.Cells.Item(94, 11).Value = 1116.4242
$ws.Cells.Item(94, 12).Value = 1684
$ws.Cells.Item(94, 13).Value = -665.4241999999999
$ws.Cells.Item(94, 14).Value = -2586

# Sheet BSM Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 56031.85
$ws.Cells.Item(134, 9).Value = 69571.06
$ws.Cells.Item(134, 10).Value = 1875
$ws.Cells.Item(134, 11).Value = 208713.18
$ws.Cells.Item(134, 12).Value = 5625
$ws.Cells.Item(134, 13).Value = -206178.18
$ws.Cells.Item(134, 14).Value = -10695

# Sheet CRP Row 20: Re-crating the Scene | Iron Spear
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 38499.75
$ws.Cells.Item(20, 10).Value = 41999.5
$ws.Cells.Item(20, 12).Value = 41999.5
$ws.Cells.Item(20, 14).Value = -42471.5

# Sheet CRP Row 30: Polearms Aplenty | Iron Spear
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(30, 8).Value = 38499.75
$ws.Cells.Item(30, 10).Value = 41999.5
$ws.Cells.Item(30, 12).Value = 41999.5
$ws.Cells.Item(30, 14).Value = -42181.5

# Sheet CRP Row 99: O Pine | Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 14709720
$ws.Cells.Item(99, 9).Value = 3086.0435
$ws.Cells.Item(99, 11).Value = 3086.0435
$ws.Cells.Item(99, 13).Value = -1588.0435

# Sheet CRP Row 126: A Better Conductor | Red Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 14709720
$ws.Cells.Item(126, 9).Value = 3086.0435
$ws.Cells.Item(126, 11).Value = 9258.130500000001
$ws.Cells.Item(126, 13).Value = -6788.130500000001

# Sheet CRP Row 128: An A-prop-riate Request | Ironwood Spear
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(128, 8).Value = 38499.75
$ws.Cells.Item(128, 10).Value = 41999.5
$ws.Cells.Item(128, 12).Value = 41999.5
$ws.Cells.Item(128, 14).Value = -51959.5

# Sheet CRP Row 134: Wood You Be Quiet | Ceiba Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1208.5385
$ws.Cells.Item(134, 9).Value = 1420.5
$ws.Cells.Item(134, 10).Value = 1144.95
$ws.Cells.Item(134, 11).Value = 4261.5
$ws.Cells.Item(134, 12).Value = 3434.85
$ws.Cells.Item(134, 13).Value = -1726.5
$ws.Cells.Item(134, 14).Value = -8504.85

# Sheet CUL Row 98: Sweet Kiss of Death | Rice Vinegar
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 349.92307
$ws.Cells.Item(98, 10).Value = 391.33334
$ws.Cells.Item(98, 12).Value = 1174.00002
$ws.Cells.Item(98, 14).Value = -4170.000019999999

# Sheet CUL Row 131: The Mountain Steeped | Tsai tou Vounou
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 728.16
$ws.Cells.Item(131, 10).Value = 728.16
$ws.Cells.Item(131, 12).Value = 2184.48
$ws.Cells.Item(131, 14).Value = -12264.48

# Sheet CUL Row 140: Sweet, Sweet Bean Juice | Mesquite Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 1564.2
$ws.Cells.Item(140, 9).Value = 1374.9445
$ws.Cells.Item(140, 11).Value = 4124.833500000001
$ws.Cells.Item(140, 13).Value = 1055.166499999999

# Sheet GSM Row 70: Sky Is the Limit | Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 11001.643
$ws.Cells.Item(70, 9).Value = 15898.375
$ws.Cells.Item(70, 10).Value = 4472.6665
$ws.Cells.Item(70, 11).Value = 15898.375
$ws.Cells.Item(70, 12).Value = 4472.6665
$ws.Cells.Item(70, 13).Value = -15628.375
$ws.Cells.Item(70, 14).Value = -5012.6665

# Sheet GSM Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 11001.643
$ws.Cells.Item(73, 9).Value = 15898.375
$ws.Cells.Item(73, 10).Value = 4472.6665
$ws.Cells.Item(73, 11).Value = 15898.375
$ws.Cells.Item(73, 12).Value = 4472.6665
$ws.Cells.Item(73, 13).Value = -14962.375
$ws.Cells.Item(73, 14).Value = -6344.6665

# Sheet GSM Row 132: On Board for Lar | Lar Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 70093.17
$ws.Cells.Item(132, 9).Value = 67920.94
$ws.Cells.Item(132, 11).Value = 203762.82
$ws.Cells.Item(132, 13).Value = -201232.82

# Sheet LTW Row 7: Tan Before the Ban | Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3213.0715
$ws.Cells.Item(7, 9).Value = 3413.7368
$ws.Cells.Item(7, 10).Value = 2789.4443
$ws.Cells.Item(7, 11).Value = 3413.7368
$ws.Cells.Item(7, 12).Value = 2789.4443
$ws.Cells.Item(7, 13).Value = -3301.7368
$ws.Cells.Item(7, 14).Value = -3013.4443

# Sheet LTW Row 61: Spelling Me Softly | Raptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3565.36
$ws.Cells.Item(61, 9).Value = 1514.5
$ws.Cells.Item(61, 11).Value = 1514.5
$ws.Cells.Item(61, 13).Value = -1312.5

# Sheet LTW Row 113: Peace in Rest | Atrociraptor Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 3565.36
$ws.Cells.Item(113, 9).Value = 1514.5
$ws.Cells.Item(113, 11).Value = 1514.5
$ws.Cells.Item(113, 13).Value = 655.5

# Sheet LTW Row 126: Battered Books | Saiga Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 3213.0715
$ws.Cells.Item(126, 9).Value = 3413.7368
$ws.Cells.Item(126, 10).Value = 2789.4443
$ws.Cells.Item(126, 11).Value = 10241.2104
$ws.Cells.Item(126, 12).Value = 8368.332900000001
$ws.Cells.Item(126, 13).Value = -7771.2104
$ws.Cells.Item(126, 14).Value = -13308.3329

# Sheet LTW Row 132: Tenets of Tanning | Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1598.8108
$ws.Cells.Item(132, 9).Value = 1023.1923
$ws.Cells.Item(132, 10).Value = 2959.3635
$ws.Cells.Item(132, 11).Value = 3069.5769
$ws.Cells.Item(132, 12).Value = 8878.0905
$ws.Cells.Item(132, 13).Value = -539.5769
$ws.Cells.Item(132, 14).Value = -13938.0905

# Sheet LTW Row 134: Freezing Fingers | Crocodileskin Fingerless Gloves of Striking
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(134, 8).Value = 52631
$ws.Cells.Item(134, 10).Value = 52631
$ws.Cells.Item(134, 12).Value = 52631
$ws.Cells.Item(134, 14).Value = -62771

# Sheet WVR Row 100: Of Great Import | Kudzu Thread
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 9).Value = 500
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 1000
$ws.Cells.Item(100, 12).ClearContents()  # L100 removed (was 1000)
$ws.Cells.Item(100, 14).Value = 0
$ws.Cells.Item(100, 13).Value = -459

# Sheet WVR Row 107: Flax Wax | Bright Linen Yarn
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2526095
$ws.Cells.Item(107, 9).Value = 978.36365
$ws.Cells.Item(107, 11).Value = 2935.09095
$ws.Cells.Item(107, 13).Value = -1015.09095

# Sheet WVR Row 122: Heavy Armoire | Dark Hempen Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1420.9512
$ws.Cells.Item(122, 9).Value = 1450
$ws.Cells.Item(122, 11).Value = 4350
$ws.Cells.Item(122, 13).Value = -1900

# Sheet WVR Row 126: A Polished Purchase | Snow Linen
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1356
$ws.Cells.Item(126, 9).Value = 1666.2727
$ws.Cells.Item(126, 10).Value = 1045.7273
$ws.Cells.Item(126, 11).Value = 4998.8181
$ws.Cells.Item(126, 12).Value = 3137.1819
$ws.Cells.Item(126, 13).Value = -2528.8181
$ws.Cells.Item(126, 14).Value = -8077.1819

# Sheet WVR Row 132: Comfy Cabins | Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 962.5714
$ws.Cells.Item(132, 9).Value = 660.1739
$ws.Cells.Item(132, 11).Value = 1980.5217
$ws.Cells.Item(132, 13).Value = 549.4783

# Sheet WVR Row 136: Weaving the Envelope | Sarcenet Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 18183202
$ws.Cells.Item(136, 9).Value = 25642230
$ws.Cells.Item(136, 10).Value = 1822.125
$ws.Cells.Item(136, 11).Value = 76926690
$ws.Cells.Item(136, 12).Value = 5466.375
$ws.Cells.Item(136, 13).Value = -76924140
$ws.Cells.Item(136, 14).Value = -10566.375
